$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calculations")
Write-Host ("B4 before: " + $ws.Range("B4").Value2)
$ws.Range("B4").Value2 = "See elec/CCaMC"
$ws.Range("B10").Value2 = "See elec/CCaMC"
Write-Host ("B4 after: " + $ws.Range("B4").Value2)
